$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the F-column formulas: change D+E to D-E for every data row (2..71)
for ($r = 2; $r -le 71; $r++) {
    $ws.Range("F$r").Formula = "=D$r-E$r"
}

# Column width changes
$ws.Columns.Item(4).ColumnWidth = 29.28515625
$ws.Columns.Item(5).ColumnWidth = 27.85546875
$ws.Columns.Item(6).ColumnWidth = 23.7109375

# Update active selection
$ws.Range("J8").Select()
